$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new data rows (136 and 137) to columns A (Marked) and B (Issue_Num)
$ws.Range("A136").Value = 0
$ws.Range("B136").Value = 240

$ws.Range("A137").Value = 1
$ws.Range("B137").Value = 243

# Scroll the view down and move the selection to the new last row
$excel.ActiveWindow.ScrollRow = 121
$ws.Range("C137").Select()
